$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180602550506592
$ws.Range("B1").Value = 1.35085117816925
$ws.Range("C1").Value = 1.128082513809204
$ws.Range("D1").Value = 4.648920059204102
$ws.Range("E1").Value = 1.525807619094849
